# slots-import-delete-fails.xlsx maintenance update:
#  - add a "Date" label to A1 (it already held the date value in B1)
#  - add a new "Help" row with a hyperlink to the CCDB conventions doc
#  - rename a handful of OPERATION picklist values to their fuller names
#  - widen column A to fit the new labels

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new "Help" row right after the Status row (row 6) ---
$ws.Rows("7").Insert()
$ws.Rows("7").RowHeight = 21.95

# Label + hyperlink cell for the new row
$ws.Range("A7").Value = "Help"
$ws.Range("A7").Style = "Normal"

$ws.Range("B7").Value = "https://ccdb.esss.lu.se/resources/help/ccdb_conventions.pdf"
$ws.Hyperlinks.Add($ws.Range("B7"), "https://ccdb.esss.lu.se/resources/help/ccdb_conventions.pdf")

# --- 2. Give the top-left date cell a label ---
$ws.Range("A1").Value = "Date"

# --- 3. Rename OPERATION values throughout column A (rows shifted down by 1 after the insert) ---
$lastRow = $ws.Cells(1, 1).SpecialCells(11).Row
for ($r = 11; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    switch ($cell.Value) {
        "CREATE RELATION" { $cell.Value = "CREATE RELATIONSHIP" }
        "DELETE RELATION" { $cell.Value = "DELETE RELATIONSHIP" }
        "INSTALL"         { $cell.Value = "INSTALL DEVICE" }
        "UNINSTALL"       { $cell.Value = "UNINSTALL DEVICE" }
    }
}

# --- 4. Update the OPERATION data validation list to match the new labels ---
$ws.Range("A11:A1048576").Validation.Delete()
$ws.Range("A11:A1048576").Validation.Add(3, 1, 1, """CREATE ENTITY,UPDATE ENTITY,DELETE ENTITY,CREATE PROPERTY,UPDATE PROPERTY,DELETE PROPERTY,CREATE RELATIONSHIP,UPDATE RELATIONSHIP,DELETE RELATIONSHIP,INSTALL DEVICE,UNINSTALL DEVICE""")
$ws.Range("A11:A1048576").Validation.InputTitle = ""

# --- 5. Widen column A so the longer labels fit ---
$ws.Columns("A:A").ColumnWidth = 20.1

# --- 6. Keep the workbook's recorded absolute path in sync with the new checkout location ---
$wb.Path_AbsPath = "C:\Projects\ess-git\ccdb\conf-core\src\test\resources\dataloader\"
